$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RUT (Chilean national ID) values to use a valid check digit,
# per commit message "ruts con digito verificador valido".
$ws.Range("C4").Value = "9876543-3"
$ws.Range("C5").Value = "20555666-4"
$ws.Range("C6").Value = "21123456-3"
$ws.Range("C7").Value = "19999888-9"
$ws.Range("C8").Value = "22333444-K"

# Move the active selection to A9 (matches saved selection state in workbook).
$ws.Range("A9").Select() | Out-Null
